# C1--C2-and-C3-PowerPoint.pptx edit
#
# 1. The table on slide 16 switches from the deck's custom table style
#    ({06803272-0A34-4E15-88DD-778CA48B251B}, "Table_0" in tableStyles.xml)
#    to the built-in PowerPoint table style {A72E9D03-5BCD-454C-B11E-7D172EB6BD67}.
#
# 2. The theme used by the slide master ("Integral") and the theme used by the
#    notes master ("Office Theme") are swapped. The only PowerPoint-object-model
#    surface that reaches theme colors is Slide.ThemeColorScheme, which maps onto
#    the slide master's theme part - so we repaint its 12 scheme colors with the
#    stock "Office" palette (what the notes-master theme originally held), which
#    is what the slide-master theme ends up holding after the swap.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 ---------------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{A72E9D03-5BCD-454C-B11E-7D172EB6BD67}", $true)
    }
}

# --- 2. Theme colour swap (slide-master theme -> stock Office palette) -----
$themeColors = $p.Slides.Item(1).ThemeColorScheme

# index : scheme slot : RGB() value for srgbClr "RRGGBB" (stored little-endian as 0x00BBGGRR)
$officePalette = @(
    0,         # 1  dk1      000000
    16777215,  # 2  lt1      FFFFFF
    6968388,   # 3  dk2      44546A
    15132391,  # 4  lt2      E7E6E6
    13998939,  # 5  accent1  5B9BD5
    3243501,   # 6  accent2  ED7D31
    10855845,  # 7  accent3  A5A5A5
    49407,     # 8  accent4  FFC000
    12874308,  # 9  accent5  4472C4
    4697456,   # 10 accent6  70AD47
    12673797,  # 11 hlink    0563C1
    7491477    # 12 folHlink 954F72
)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officePalette[$i - 1]
}
